{"js": "// Helper: find the (unique) paragraph that contains `needle` and replace\n// its OOXML with `paragraphOoxml` (a single <w:p>...</w:p> fragment).\nasync function replaceParagraphByText(context, needle, paragraphOoxml) {\n  const results = context.document.body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"paragraph not found for: \" + needle);\n  }\n  const range = results.items[0];\n  const paras = range.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const para = paras.items[0];\n\n  const pkg =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + paragraphOoxml + '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n\n  para.insertOoxml(pkg, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst GREEN_RPR =\n  '<w:rPr><w:color w:val=\"538135\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr>';\nconst GREEN_STRIKE_RPR =\n  '<w:rPr><w:strike/><w:color w:val=\"538135\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr>';\n\n// 1) \"...i.e., being on 11 and rolling a 4 will result in a new position of 2.\"\n//    -> \"...of 3.\" (split the trailing sentence into its own 3 runs: text/\"3\"/\".\")\nawait replaceParagraphByText(\n  context,\n  \"being on 11 and rolling a 4 will result in a new position of 2.\",\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    GREEN_RPR +\n    '</w:pPr>' +\n    '<w:r>' + GREEN_RPR + '<w:t xml:space=\"preserve\">Verify that </w:t></w:r>' +\n    '<w:r>' + GREEN_RPR + '<w:t xml:space=\"preserve\">(in penalty box) </w:t></w:r>' +\n    '<w:r>' + GREEN_RPR + '<w:t xml:space=\"preserve\">rolling more than an 11 will result in the board resetting and starting from the beginning with the additional positional </w:t></w:r>' +\n    '<w:r>' + GREEN_RPR + '<w:t xml:space=\"preserve\">increases from the roll. </w:t></w:r>' +\n    '<w:r>' + GREEN_RPR + '<w:t>i.e.,</w:t></w:r>' +\n    '<w:r>' + GREEN_RPR + '<w:t xml:space=\"preserve\"> being on 11 and rolling a 4 will result in a new position of </w:t></w:r>' +\n    '<w:r>' + GREEN_RPR + '<w:t>3</w:t></w:r>' +\n    '<w:r>' + GREEN_RPR + '<w:t>.</w:t></w:r>' +\n    '</w:p>'\n);\n\n// 2) \"Verify that (out of the penalty box) a roll will be added...\" gains\n//    strikethrough + green formatting.\nawait replaceParagraphByText(\n  context,\n  \"Verify that (out of the penalty box) a roll will be added\",\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    GREEN_STRIKE_RPR +\n    '</w:pPr>' +\n    '<w:r>' + GREEN_STRIKE_RPR +\n    '<w:t>Verify that (out of the penalty box) a roll will be added to the current position of the player (+=)</w:t></w:r>' +\n    '</w:p>'\n);\n\n// 3) \"Verify that (out of the penalty box) rolling more than an 11...\" gains\n//    strikethrough + green formatting on all 3 runs (text itself is unchanged,\n//    it already ends in \"...of 3.\").\nawait replaceParagraphByText(\n  context,\n  \"Verify that (out of the penalty box) rolling more than an 11\",\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    GREEN_STRIKE_RPR +\n    '</w:pPr>' +\n    '<w:r>' + GREEN_STRIKE_RPR +\n    '<w:t xml:space=\"preserve\">Verify that (out of the penalty box) rolling more than an 11 will result in the board resetting and starting from the beginning with the additional positional increases from the roll. i.e., being on 11 and rolling a 4 will result in a new position of </w:t></w:r>' +\n    '<w:r>' + GREEN_STRIKE_RPR + '<w:t>3</w:t></w:r>' +\n    '<w:r>' + GREEN_STRIKE_RPR + '<w:t>.</w:t></w:r>' +\n    '</w:p>'\n);\n\n// 4) \"Verify that a person's position, is the culmination...\" gains\n//    strikethrough + green formatting.\nawait replaceParagraphByText(\n  context,\n  \"culmination of their total rolls\",\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr>' +\n    GREEN_STRIKE_RPR +\n    '</w:pPr>' +\n    '<w:r>' + GREEN_STRIKE_RPR +\n    '<w:t>Verify that a person\\u2019s position, is the culmination of their total rolls. Rolling 5 and then rolling 2, will have them at position 7.</w:t></w:r>' +\n    '</w:p>'\n);\n\n// 5) \"Verify that all position on the board correspond to the correct\n//    category:\" -> \"...all positions on the board correspond to the correct\n//    categories:\" (plain text, no color/strike formatting).\nawait replaceParagraphByText(\n  context,\n  \"position on the board\",\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Verify that </w:t></w:r>' +\n    '<w:r><w:t>all</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> position</w:t></w:r>' +\n    '<w:r><w:t>s</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> on the board </w:t></w:r>' +\n    '<w:r><w:t>correspond</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:t>to</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> the correct categor</w:t></w:r>' +\n    '<w:r><w:t>ies</w:t></w:r>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n    '</w:p>'\n);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($needle) {\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text.Contains($needle)) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Replace-ParagraphXml($needle, $innerParagraphXml) {\n    $idx = Find-ParagraphIndex($needle)\n    if ($idx -eq -1) {\n        throw \"paragraph not found for: $needle\"\n    }\n    $p = $d.Paragraphs.Item($idx)\n    $r = $p.Range\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $innerParagraphXml + '</w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n    [void]$r.InsertXML($xml)\n}\n\n$greenRpr = '<w:rPr><w:color w:val=\"538135\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr>'\n$greenStrikeRpr = '<w:rPr><w:strike/><w:color w:val=\"538135\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr>'\n\n# 1) \"...i.e., being on 11 and rolling a 4 will result in a new position of 2.\"\n#    -> \"...of 3.\" (trailing sentence ends up split into 3 runs: text/\"3\"/\".\")\n$para17 = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    $greenRpr + '</w:pPr>' +\n    '<w:r>' + $greenRpr + '<w:t xml:space=\"preserve\">Verify that </w:t></w:r>' +\n    '<w:r>' + $greenRpr + '<w:t xml:space=\"preserve\">(in penalty box) </w:t></w:r>' +\n    '<w:r>' + $greenRpr + '<w:t xml:space=\"preserve\">rolling more than an 11 will result in the board resetting and starting from the beginning with the additional positional </w:t></w:r>' +\n    '<w:r>' + $greenRpr + '<w:t xml:space=\"preserve\">increases from the roll. </w:t></w:r>' +\n    '<w:r>' + $greenRpr + '<w:t>i.e.,</w:t></w:r>' +\n    '<w:r>' + $greenRpr + '<w:t xml:space=\"preserve\"> being on 11 and rolling a 4 will result in a new position of </w:t></w:r>' +\n    '<w:r>' + $greenRpr + '<w:t>3</w:t></w:r>' +\n    '<w:r>' + $greenRpr + '<w:t>.</w:t></w:r>' +\n    '</w:p>'\nReplace-ParagraphXml \"being on 11 and rolling a 4 will result in a new position of 2.\" $para17\n\n# 2) \"Verify that (out of the penalty box) a roll will be added...\" gains\n#    strikethrough + green formatting.\n$para21 = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    $greenStrikeRpr + '</w:pPr>' +\n    '<w:r>' + $greenStrikeRpr +\n    '<w:t>Verify that (out of the penalty box) a roll will be added to the current position of the player (+=)</w:t></w:r>' +\n    '</w:p>'\nReplace-ParagraphXml \"Verify that (out of the penalty box) a roll will be added\" $para21\n\n# 3) \"Verify that (out of the penalty box) rolling more than an 11...\" gains\n#    strikethrough + green formatting on all 3 runs (text is unchanged, it\n#    already ends in \"...of 3.\").\n$para22 = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    $greenStrikeRpr + '</w:pPr>' +\n    '<w:r>' + $greenStrikeRpr +\n    '<w:t xml:space=\"preserve\">Verify that (out of the penalty box) rolling more than an 11 will result in the board resetting and starting from the beginning with the additional positional increases from the roll. i.e., being on 11 and rolling a 4 will result in a new position of </w:t></w:r>' +\n    '<w:r>' + $greenStrikeRpr + '<w:t>3</w:t></w:r>' +\n    '<w:r>' + $greenStrikeRpr + '<w:t>.</w:t></w:r>' +\n    '</w:p>'\nReplace-ParagraphXml \"Verify that (out of the penalty box) rolling more than an 11\" $para22\n\n# 4) \"Verify that a person's position, is the culmination...\" gains\n#    strikethrough + green formatting.\n$para24 = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr>' +\n    $greenStrikeRpr + '</w:pPr>' +\n    '<w:r>' + $greenStrikeRpr +\n    '<w:t>Verify that a person' + [char]0x2019 + 's position, is the culmination of their total rolls. Rolling 5 and then rolling 2, will have them at position 7.</w:t></w:r>' +\n    '</w:p>'\nReplace-ParagraphXml \"culmination of their total rolls\" $para24\n\n# 5) \"Verify that all position on the board correspond to the correct\n#    category:\" -> \"...all positions on the board correspond to the correct\n#    categories:\" (plain text, no color/strike formatting).\n$para25 = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Verify that </w:t></w:r>' +\n    '<w:r><w:t>all</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> position</w:t></w:r>' +\n    '<w:r><w:t>s</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> on the board </w:t></w:r>' +\n    '<w:r><w:t>correspond</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:t>to</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> the correct categor</w:t></w:r>' +\n    '<w:r><w:t>ies</w:t></w:r>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n    '</w:p>'\nReplace-ParagraphXml \"position on the board\" $para25\n\nWrite-Output \"done\"\n"}
